$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub6 = [string][char]0x2086   # U+2086 SUBSCRIPT SIX, used in one price value

# Updated crypto price/volume data: cell reference -> new text value
$updates = @(
    @{ Cell = "D2"; Value = '67.836.04' }
    @{ Cell = "E2"; Value = '  +1.23%  ' }
    @{ Cell = "D3"; Value = '2.619.10' }
    @{ Cell = "E3"; Value = '  +0.37%  ' }
    @{ Cell = "E4"; Value = '  -0.14%  ' }
    @{ Cell = "D5"; Value = '595.46' }
    @{ Cell = "E5"; Value = '  +0.53%  ' }
    @{ Cell = "D6"; Value = '152.88' }
    @{ Cell = "E6"; Value = '  +0.37%  ' }
    @{ Cell = "E7"; Value = '  +0.01%  ' }
    @{ Cell = "E8"; Value = '  -1.89%  ' }
    @{ Cell = "D9"; Value = '2.619.26' }
    @{ Cell = "E9"; Value = '  +0.40%  ' }
    @{ Cell = "D10"; Value = '0.133' }
    @{ Cell = "E10"; Value = '  +9.48%  ' }
    @{ Cell = "E11"; Value = '  -0.56%  ' }
    @{ Cell = "D12"; Value = '5.23' }
    @{ Cell = "E12"; Value = '  +0.96%  ' }
    @{ Cell = "D13"; Value = '0.346' }
    @{ Cell = "E13"; Value = '  +0.09%  ' }
    @{ Cell = "D14"; Value = '27.57' }
    @{ Cell = "E14"; Value = '  +0.30%  ' }
    @{ Cell = "D15"; Value = '0.0000187' }
    @{ Cell = "E15"; Value = '  +4.37%  ' }
    @{ Cell = "D16"; Value = '3.095.32' }
    @{ Cell = "E16"; Value = '  +0.39%  ' }
    @{ Cell = "D17"; Value = '67.667.41' }
    @{ Cell = "E17"; Value = '  +1.13%  ' }
    @{ Cell = "D18"; Value = '2.621.26' }
    @{ Cell = "E18"; Value = '  -0.12%  ' }
    @{ Cell = "D19"; Value = '11.36' }
    @{ Cell = "E19"; Value = '  +3.11%  ' }
    @{ Cell = "D20"; Value = '367.40' }
    @{ Cell = "E20"; Value = '  +0.56%  ' }
    @{ Cell = "D21"; Value = '7.39' }
    @{ Cell = "E21"; Value = '  +0.70%  ' }
    @{ Cell = "D22"; Value = '4.22' }
    @{ Cell = "E22"; Value = '  -1.68%  ' }
    @{ Cell = "D23"; Value = '4.77' }
    @{ Cell = "E23"; Value = '  -1.28%  ' }
    @{ Cell = "D24"; Value = '2.06' }
    @{ Cell = "E24"; Value = '  +0.83%  ' }
    @{ Cell = "D25"; Value = '72.66' }
    @{ Cell = "E25"; Value = '  +9.57%  ' }
    @{ Cell = "D26"; Value = '1.00' }
    @{ Cell = "E26"; Value = '  -0.10%  ' }
    @{ Cell = "D27"; Value = '9.92' }
    @{ Cell = "E27"; Value = '  -1.44%  ' }
    @{ Cell = "D28"; Value = '2.750.64' }
    @{ Cell = "E28"; Value = '  +0.41%  ' }
    @{ Cell = "D29"; Value = '0.0000103' }
    @{ Cell = "E29"; Value = '  +2.43%  ' }
    @{ Cell = "D30"; Value = '1.00' }
    @{ Cell = "E30"; Value = '  +0.29%  ' }
    @{ Cell = "D31"; Value = '568.80' }
    @{ Cell = "E31"; Value = '  -2.67%  ' }
    @{ Cell = "D32"; Value = '7.86' }
    @{ Cell = "E32"; Value = '  +2.10%  ' }
    @{ Cell = "D33"; Value = '1.39' }
    @{ Cell = "E33"; Value = '  +0.90%  ' }
    @{ Cell = "D34"; Value = '1.83' }
    @{ Cell = "E34"; Value = '  +1.37%  ' }
    @{ Cell = "D35"; Value = '1.00' }
    @{ Cell = "E35"; Value = '  +0.04%  ' }
    @{ Cell = "D36"; Value = '0.127' }
    @{ Cell = "E36"; Value = '  +3.17%  ' }
    @{ Cell = "E37"; Value = '  +1.09%  ' }
    @{ Cell = "D38"; Value = '161.00' }
    @{ Cell = "E38"; Value = '  +3.81%  ' }
    @{ Cell = "D39"; Value = '19.07' }
    @{ Cell = "E39"; Value = '  +1.03%  ' }
    @{ Cell = "D40"; Value = '1.87' }
    @{ Cell = "E40"; Value = '  +4.00%  ' }
    @{ Cell = "D41"; Value = '0.366' }
    @{ Cell = "E41"; Value = '  +0.19%  ' }
    @{ Cell = "D42"; Value = '5.28' }
    @{ Cell = "E42"; Value = '  +1.40%  ' }
    @{ Cell = "D43"; Value = '2.64' }
    @{ Cell = "E43"; Value = '  +2.55%  ' }
    @{ Cell = "D44"; Value = ('0.0' + $sub6 + '0327') }
    @{ Cell = "E44"; Value = '  +12.14%  ' }
    @{ Cell = "D45"; Value = '17.31' }
    @{ Cell = "E45"; Value = '  +4.08%  ' }
    @{ Cell = "D46"; Value = '1.00' }
    @{ Cell = "E46"; Value = '  +0.09%  ' }
    @{ Cell = "D47"; Value = '40.16' }
    @{ Cell = "E47"; Value = '  -1.51%  ' }
    @{ Cell = "D48"; Value = '154.79' }
    @{ Cell = "E48"; Value = '  +0.19%  ' }
    @{ Cell = "D49"; Value = '3.66' }
    @{ Cell = "E49"; Value = '  -1.23%  ' }
    @{ Cell = "D50"; Value = '21.70' }
    @{ Cell = "E50"; Value = '  +1.04%  ' }
    @{ Cell = "D51"; Value = '1.69' }
    @{ Cell = "E51"; Value = '  -0.12%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Prefix with an apostrophe to force text, so numeric-looking strings
    # (e.g. '1.00', '67.836.04') are not reinterpreted as numbers/dates.
    $cell.Value = "'" + [string]$u.Value
    # Reset style back to Normal so no stray text-format style gets applied
    $cell.Style = "Normal"
}
